$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 for the new facility (Chemours Chambers Works / Deepwater),
# shifting the existing rows 4-10 down to rows 5-11.
$ws.Rows.Item(4).Insert()

# Write final values for every data row (2-11), covering columns A-S.
# (Blank GHG_co2e cells for rows 6 and 10 are left untouched -- Insert() already
# carried their blank state forward from the original rows 5 and 9.)

# Row 2: Daikin America Inc.
$ws.Cells.Item(2, 1).Value = 'Daikin America Inc.'
$ws.Cells.Item(2, 2).Value = 'Decatur'
$ws.Cells.Item(2, 3).Value = 96067
$ws.Cells.Item(2, 4).Value = 37
$ws.Cells.Item(2, 5).Value = 284.138128274564
$ws.Cells.Item(2, 6).Value = 53051
$ws.Cells.Item(2, 7).Value = 186.708486897389
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0.22
$ws.Cells.Item(2, 10).Value = 33838
$ws.Cells.Item(2, 11).Value = 13456
$ws.Cells.Item(2, 12).Value = 495
$ws.Cells.Item(2, 13).Value = 156
$ws.Cells.Item(2, 14).Value = 7531
$ws.Cells.Item(2, 15).Value = 47.6078918918919
$ws.Cells.Item(2, 16).Value = 7.4115680192577
$ws.Cells.Item(2, 17).Value = 13.1056293233406
$ws.Cells.Item(2, 18).Value = 39.4594594594595
$ws.Cells.Item(2, 19).Value = 0.486486486486487

# Row 3: Chemours El Dorado
$ws.Cells.Item(3, 1).Value = 'Chemours El Dorado'
$ws.Cells.Item(3, 2).Value = 'El Dorado'
$ws.Cells.Item(3, 3).Value = 66990
$ws.Cells.Item(3, 4).Value = 6
$ws.Cells.Item(3, 5).Value = 422.192389752563
$ws.Cells.Item(3, 6).Value = 8797
$ws.Cells.Item(3, 7).Value = 20.8364722186388
$ws.Cells.Item(3, 8).Value = 1
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 7211
$ws.Cells.Item(3, 11).Value = 1335
$ws.Cells.Item(3, 12).Value = 173
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 14).Value = 298
$ws.Cells.Item(3, 15).Value = 55.6223333333333
$ws.Cells.Item(3, 16).Value = 4.00529428733787
$ws.Cells.Item(3, 17).Value = 11.1294509931195
$ws.Cells.Item(3, 18).Value = 50
$ws.Cells.Item(3, 19).Value = 0.566666666666667

# Row 4: Chemours Chambers Works
$ws.Cells.Item(4, 1).Value = 'Chemours Chambers Works'
$ws.Cells.Item(4, 2).Value = 'Deepwater'
$ws.Cells.Item(4, 3).Value = 2619
$ws.Cells.Item(4, 4).Value = 71
$ws.Cells.Item(4, 5).Value = 140.908502645324
$ws.Cells.Item(4, 6).Value = 86712
$ws.Cells.Item(4, 7).Value = 615.378052935953
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0.23
$ws.Cells.Item(4, 10).Value = 46412
$ws.Cells.Item(4, 11).Value = 32532
$ws.Cells.Item(4, 12).Value = 263
$ws.Cells.Item(4, 13).Value = 1857
$ws.Cells.Item(4, 14).Value = 11099
$ws.Cells.Item(4, 15).Value = 53.0347246376812
$ws.Cells.Item(4, 16).Value = 9.19041479039809
$ws.Cells.Item(4, 17).Value = 10.1818382334019
$ws.Cells.Item(4, 18).Value = 39
$ws.Cells.Item(4, 19).Value = 0.365714285714286

# Row 5: Chemours Louisville Works
$ws.Cells.Item(5, 1).Value = 'Chemours Louisville Works'
$ws.Cells.Item(5, 2).Value = 'Louisville'
$ws.Cells.Item(5, 3).Value = 3707770
$ws.Cells.Item(5, 4).Value = 188
$ws.Cells.Item(5, 5).Value = 128.555980630485
$ws.Cells.Item(5, 6).Value = 203331
$ws.Cells.Item(5, 7).Value = 1581.65337001664
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0.016
$ws.Cells.Item(5, 10).Value = 102835
$ws.Cells.Item(5, 11).Value = 88097
$ws.Cells.Item(5, 12).Value = 487
$ws.Cells.Item(5, 13).Value = 3577
$ws.Cells.Item(5, 14).Value = 9205
$ws.Cells.Item(5, 15).Value = 36.6684770114943
$ws.Cells.Item(5, 16).Value = 12.486207688811
$ws.Cells.Item(5, 17).Value = 14.6084983177459
$ws.Cells.Item(5, 18).Value = 30.3208556149733
$ws.Cells.Item(5, 19).Value = 0.429411764705882

# Row 6: Iofina Chemical Inc.
$ws.Cells.Item(6, 1).Value = 'Iofina Chemical Inc.'
$ws.Cells.Item(6, 2).Value = 'Covington'
$ws.Cells.Item(6, 4).Value = 110
$ws.Cells.Item(6, 5).Value = 119.931178899633
$ws.Cells.Item(6, 6).Value = 175554
$ws.Cells.Item(6, 7).Value = 1463.7894966989
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0.091
$ws.Cells.Item(6, 10).Value = 158873
$ws.Cells.Item(6, 11).Value = 7479
$ws.Cells.Item(6, 12).Value = 278
$ws.Cells.Item(6, 13).Value = 2452
$ws.Cells.Item(6, 14).Value = 5971
$ws.Cells.Item(6, 15).Value = 64.9617889908257
$ws.Cells.Item(6, 16).Value = 5.81435262743146
$ws.Cells.Item(6, 17).Value = 6.80743447224533
$ws.Cells.Item(6, 18).Value = 30
$ws.Cells.Item(6, 19).Value = 0.393636363636364

# Row 7: ARKEMA, INC.
$ws.Cells.Item(7, 1).Value = 'ARKEMA, INC.'
$ws.Cells.Item(7, 2).Value = 'Calvert City'
$ws.Cells.Item(7, 3).Value = 843010
$ws.Cells.Item(7, 4).Value = 11
$ws.Cells.Item(7, 5).Value = 200.196602297996
$ws.Cells.Item(7, 6).Value = 13612
$ws.Cells.Item(7, 7).Value = 67.9931619405725
$ws.Cells.Item(7, 8).Value = 1
$ws.Cells.Item(7, 9).Value = 1
$ws.Cells.Item(7, 10).Value = 13385
$ws.Cells.Item(7, 11).Value = 77
$ws.Cells.Item(7, 12).Value = 23
$ws.Cells.Item(7, 13).Value = 10
$ws.Cells.Item(7, 14).Value = 384
$ws.Cells.Item(7, 15).Value = 57.363
$ws.Cells.Item(7, 16).Value = 6.63943958765135
$ws.Cells.Item(7, 17).Value = 4.02491570842668
$ws.Cells.Item(7, 18).Value = 33.6363636363636
$ws.Cells.Item(7, 19).Value = 1.26363636363636

# Row 8: Honeywell International - Geismar Complex
$ws.Cells.Item(8, 1).Value = 'Honeywell International - Geismar Complex'
$ws.Cells.Item(8, 2).Value = 'Geismar'
$ws.Cells.Item(8, 3).Value = 413584
$ws.Cells.Item(8, 4).Value = 14
$ws.Cells.Item(8, 5).Value = 228.365522556404
$ws.Cells.Item(8, 6).Value = 38167
$ws.Cells.Item(8, 7).Value = 167.131183257197
$ws.Cells.Item(8, 8).Value = 1
$ws.Cells.Item(8, 9).Value = 0.93
$ws.Cells.Item(8, 10).Value = 23475
$ws.Cells.Item(8, 11).Value = 13551
$ws.Cells.Item(8, 12).Value = 26
$ws.Cells.Item(8, 13).Value = 199
$ws.Cells.Item(8, 14).Value = 1091
$ws.Cells.Item(8, 15).Value = 67.4632307692308
$ws.Cells.Item(8, 16).Value = 6.08715442194185
$ws.Cells.Item(8, 17).Value = 3.24084657418216
$ws.Cells.Item(8, 18).Value = 115.714285714286
$ws.Cells.Item(8, 19).Value = 0.55

# Row 9: Mexichem Fluor Inc.
$ws.Cells.Item(9, 1).Value = 'Mexichem Fluor Inc.'
$ws.Cells.Item(9, 2).Value = 'Saint Gabriel'
$ws.Cells.Item(9, 3).Value = 18331
$ws.Cells.Item(9, 4).Value = 11
$ws.Cells.Item(9, 5).Value = 206.749916306705
$ws.Cells.Item(9, 6).Value = 28984
$ws.Cells.Item(9, 7).Value = 140.188690364466
$ws.Cells.Item(9, 8).Value = 1
$ws.Cells.Item(9, 9).Value = 1
$ws.Cells.Item(9, 10).Value = 16679
$ws.Cells.Item(9, 11).Value = 11671
$ws.Cells.Item(9, 12).Value = 26
$ws.Cells.Item(9, 13).Value = 29
$ws.Cells.Item(9, 14).Value = 722
$ws.Cells.Item(9, 15).Value = 65.1825
$ws.Cells.Item(9, 16).Value = 5.5361169432872
$ws.Cells.Item(9, 17).Value = 3.3293309509454
$ws.Cells.Item(9, 18).Value = 124.545454545455
$ws.Cells.Item(9, 19).Value = 0.563636363636364

# Row 10: Islechem LLC
$ws.Cells.Item(10, 1).Value = 'Islechem LLC'
$ws.Cells.Item(10, 2).Value = 'Grand Island'
$ws.Cells.Item(10, 4).Value = 83
$ws.Cells.Item(10, 5).Value = 81.1518555780108
$ws.Cells.Item(10, 6).Value = 86046
$ws.Cells.Item(10, 7).Value = 1060.30847214928
$ws.Cells.Item(10, 8).Value = 1
$ws.Cells.Item(10, 9).Value = 0.096
$ws.Cells.Item(10, 10).Value = 69335
$ws.Cells.Item(10, 11).Value = 10506
$ws.Cells.Item(10, 12).Value = 801
$ws.Cells.Item(10, 13).Value = 1366
$ws.Cells.Item(10, 14).Value = 3292
$ws.Cells.Item(10, 15).Value = 49.6234615384615
$ws.Cells.Item(10, 16).Value = 9.98462726344703
$ws.Cells.Item(10, 17).Value = 13.1979319795331
$ws.Cells.Item(10, 18).Value = 20.4819277108434
$ws.Cells.Item(10, 19).Value = 0.255421686746988

# Row 11: Chemours - Corpus Christi Plant
$ws.Cells.Item(11, 1).Value = 'Chemours - Corpus Christi Plant'
$ws.Cells.Item(11, 2).Value = 'Gregory'
$ws.Cells.Item(11, 3).Value = 17240
$ws.Cells.Item(11, 4).Value = 21
$ws.Cells.Item(11, 5).Value = 334.263994048958
$ws.Cells.Item(11, 6).Value = 35609
$ws.Cells.Item(11, 7).Value = 106.529571338708
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0.71
$ws.Cells.Item(11, 10).Value = 32651
$ws.Cells.Item(11, 11).Value = 768
$ws.Cells.Item(11, 12).Value = 192
$ws.Cells.Item(11, 13).Value = 303
$ws.Cells.Item(11, 14).Value = 15824
$ws.Cells.Item(11, 15).Value = 68.7501578947368
$ws.Cells.Item(11, 16).Value = 4.22768148216296
$ws.Cells.Item(11, 17).Value = 3.02033713148773
$ws.Cells.Item(11, 18).Value = 20
$ws.Cells.Item(11, 19).Value = 0.20952380952381
